$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-11 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-12 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("73×11=803", $true, $false, $false, $false, $false, $true, 1, $false, "52×43=2236", 2) | Out-Null
$d.Content.Find.Execute("55×60=3300", $true, $false, $false, $false, $false, $true, 1, $false, "59×53=3127", 2) | Out-Null
$d.Content.Find.Execute("58×61=3538", $true, $false, $false, $false, $false, $true, 1, $false, "32×97=3104", 2) | Out-Null
$d.Content.Find.Execute("84×49=4116", $true, $false, $false, $false, $false, $true, 1, $false, "78×22=1716", 2) | Out-Null
$d.Content.Find.Execute("94×39=3666", $true, $false, $false, $false, $false, $true, 1, $false, "44×46=2024", 2) | Out-Null
$d.Content.Find.Execute("95×33=3135", $true, $false, $false, $false, $false, $true, 1, $false, "42×24=1008", 2) | Out-Null
$d.Content.Find.Execute("41×26=1066", $true, $false, $false, $false, $false, $true, 1, $false, "93×36=3348", 2) | Out-Null
$d.Content.Find.Execute("80×39=3120", $true, $false, $false, $false, $false, $true, 1, $false, "60×51=3060", 2) | Out-Null
$d.Content.Find.Execute("40×23=920", $true, $false, $false, $false, $false, $true, 1, $false, "63×95=5985", 2) | Out-Null
$d.Content.Find.Execute("83×65=5395", $true, $false, $false, $false, $false, $true, 1, $false, "72×83=5976", 2) | Out-Null
$d.Content.Find.Execute("15×33=495", $true, $false, $false, $false, $false, $true, 1, $false, "66×30=1980", 2) | Out-Null
$d.Content.Find.Execute("71×74=5254", $true, $false, $false, $false, $false, $true, 1, $false, "94×12=1128", 2) | Out-Null
$d.Content.Find.Execute("66×84=5544", $true, $false, $false, $false, $false, $true, 1, $false, "60×18=1080", 2) | Out-Null
$d.Content.Find.Execute("53×79=4187", $true, $false, $false, $false, $false, $true, 1, $false, "66×55=3630", 2) | Out-Null
$d.Content.Find.Execute("28×14=392", $true, $false, $false, $false, $false, $true, 1, $false, "40×38=1520", 2) | Out-Null
$d.Content.Find.Execute("86×64=5504", $true, $false, $false, $false, $false, $true, 1, $false, "35×91=3185", 2) | Out-Null
$d.Content.Find.Execute("59×21=1239", $true, $false, $false, $false, $false, $true, 1, $false, "54×53=2862", 2) | Out-Null
$d.Content.Find.Execute("77×99=7623", $true, $false, $false, $false, $false, $true, 1, $false, "14×34=476", 2) | Out-Null
$d.Content.Find.Execute("79×15=1185", $true, $false, $false, $false, $false, $true, 1, $false, "14×47=658", 2) | Out-Null
$d.Content.Find.Execute("41×61=2501", $true, $false, $false, $false, $false, $true, 1, $false, "94×74=6956", 2) | Out-Null
$d.Content.Find.Execute("15×97=1455", $true, $false, $false, $false, $false, $true, 1, $false, "87×64=5568", 2) | Out-Null
$d.Content.Find.Execute("40×92=3680", $true, $false, $false, $false, $false, $true, 1, $false, "54×71=3834", 2) | Out-Null
$d.Content.Find.Execute("53×17=901", $true, $false, $false, $false, $false, $true, 1, $false, "76×42=3192", 2) | Out-Null
$d.Content.Find.Execute("29×70=2030", $true, $false, $false, $false, $false, $true, 1, $false, "16×59=944", 2) | Out-Null
$d.Content.Find.Execute("87×88=7656", $true, $false, $false, $false, $false, $true, 1, $false, "40×34=1360", 2) | Out-Null
